# Atualiza os dados de faturamento diário das lojas Bibi.
# - Reordena os rótulos das lojas nas linhas 2-4 (Mundi, Vieiralves, Manauara)
#   mantendo os valores de faturamento já existentes junto com cada loja.
# - Adiciona novos valores de faturamento (coluna E / dia 4) para as linhas 2-6.
# - Atualiza os totais (coluna AG) de cada linha para refletir os novos valores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha 2: Bibi Cell Mundi
$ws.Range("A2").Value = "Bibi Cell Mundi"
$ws.Range("B2").Value = 8258
$ws.Range("C2").Value = 2278
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 15207.62
$ws.Range("AG2").Value = 25743.62

# Linha 3: Bibi Cell Vieiralves
$ws.Range("A3").Value = "Bibi Cell Vieiralves"
$ws.Range("B3").Value = 8802
$ws.Range("C3").Value = 7274
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3319.4
$ws.Range("AG3").Value = 19395.4

# Linha 4: Bibi Cell Manauara
$ws.Range("A4").Value = "Bibi Cell Manauara"
$ws.Range("B4").Value = 2469.75
$ws.Range("C4").Value = 5177
$ws.Range("D4").Value = 3030
$ws.Range("E4").Value = 4202
$ws.Range("AG4").Value = 14878.75

# Linha 5: Bibi Cell Ponta Negra (rótulo inalterado)
$ws.Range("E5").Value = 2180
$ws.Range("AG5").Value = 12398.55

# Linha 6: total (rótulo inalterado)
$ws.Range("E6").Value = 24909.02
$ws.Range("AG6").Value = 72416.32000000001
